# Update prediction results: new training/testing data produced updated
# confidence values (column D) and, for some rows, a changed predicted
# status (column C) from "Aprobado" to "Desaprobado".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new CONFIANZA_APROBACION (%) value (column D)
$confianza = @{
    2  = 80.86
    3  = 83.25
    4  = 45.68
    5  = 26.73
    6  = 51.8
    7  = 58.38
    8  = 54.03
    9  = 50.26
    10 = 62.45
    11 = 45.71
    12 = 57.17
    13 = 42.94
    14 = 55.86
    15 = 50.87
    16 = 50.29
    17 = 23.84
    18 = 45.22
    19 = 45.16
    20 = 51.89
    21 = 72.1
    22 = 55.86
    23 = 73.56
    24 = 72.54
    25 = 70.83
    26 = 39.39
    27 = 29.47
}

# Rows whose PREDICCION_ESTADO (column C) flips from Aprobado to Desaprobado
$estadoDesaprobado = @(4, 11, 13, 18, 19)

foreach ($row in $confianza.Keys) {
    $ws.Cells.Item($row, 4).Value = $confianza[$row]
}

foreach ($row in $estadoDesaprobado) {
    $ws.Cells.Item($row, 3).Value = "Desaprobado"
}
